$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 04:52"

# --- Updated case counts for several countries (new data snapshot) ---
# Bolivia (row 32)
$ws.Cells.Item(32,2).Value = 89999
$ws.Cells.Item(32,3).Value = 944
$ws.Cells.Item(32,4).Value = 29808
$ws.Cells.Item(32,5).Value = 56551
$ws.Cells.Item(32,7).Value = 53
$ws.Cells.Item(32,8).Value = 3640

# Belgica (row 40)
$ws.Cells.Item(40,2).Value = 74152
$ws.Cells.Item(40,3).Value = 751
$ws.Cells.Item(40,4).Value = 17780
$ws.Cells.Item(40,5).Value = 46500
$ws.Cells.Item(40,7).Value = 2
$ws.Cells.Item(40,8).Value = 9872

# Honduras (row 50)
$ws.Cells.Item(50,2).Value = 47454
$ws.Cells.Item(50,3).Value = 481
$ws.Cells.Item(50,4).Value = 6597
$ws.Cells.Item(50,5).Value = 39362
$ws.Cells.Item(50,7).Value = 19
$ws.Cells.Item(50,8).Value = 1495

# Hong Kong (row 111)
$ws.Cells.Item(111,5).Value = 1181
$ws.Cells.Item(111,7).Value = 1
$ws.Cells.Item(111,8).Value = 52

# Vietnam (row 158)
$ws.Cells.Item(158,5).Value = 433
$ws.Cells.Item(158,7).Value = 2
$ws.Cells.Item(158,8).Value = 13

# Camboya (row 179)
$ws.Cells.Item(179,2).Value = 251
$ws.Cells.Item(179,3).Value = 3
$ws.Cells.Item(179,4).Value = 219
$ws.Cells.Item(179,5).Value = 32

# Monaco (row 189)
$ws.Cells.Item(189,2).Value = 133
$ws.Cells.Item(189,4).Value = 113
$ws.Cells.Item(189,5).Value = 16

# San Vicente y las Granadinas (row 196)
$ws.Cells.Item(196,2).Value = 57
$ws.Cells.Item(196,3).Value = 1
$ws.Cells.Item(196,4).Value = 52
$ws.Cells.Item(196,5).Value = 5

# --- Swap the Islas Malvinas / Montserrat rows (213 <-> 214) ---
# Row 213 held Islas Malvinas, row 214 held Montserrat; the refreshed
# source data now lists Montserrat first, Islas Malvinas second, each
# carrying its own row of figures along with it.
$ws.Cells.Item(213,1).Value = "Montserrat"
$ws.Cells.Item(213,4).Value = 12
$ws.Cells.Item(213,8).Value = 1

$ws.Cells.Item(214,1).Value = "Islas Malvinas"
$ws.Cells.Item(214,4).Value = 13
$ws.Cells.Item(214,8).Value = 0
